# Update InPrice (L) / OutPrice (M) values from refreshed masterdata feed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

$ws.Range("L2:L19").Value = 528.25
$ws.Range("L20:L29").Value = 452.75
$ws.Range("L30:L34").Value = 528.25
$ws.Range("L40:L43").Value = 95
$ws.Range("L48:L50").Value = 679
$ws.Range("L58:L60").Value = 566
$ws.Range("L68:L71").Value = 622.5
$ws.Range("L81:L84").Value = 566
$ws.Range("L89:L96").Value = 283
$ws.Range("L115:L118").Value = 226.25
$ws.Range("L121:L126").Value = 792.25
$ws.Range("L133:L135").Value = 1132
$ws.Range("L140:L142").Value = 566
$ws.Range("L146:L148").Value = 1075.25
$ws.Range("L155:L158").Value = 339.5
$ws.Range("L167:L169").Value = 452.75
$ws.Range("L174:L180").Value = 792.25

$ws.Range("M115:M118").Value = 599

# Restore the active selection left by the editor
$ws.Range("L13").Select()
